$d = $word.ActiveDocument

# 1. Change the port number 8080 -> 8090 in the two "访问" usage paragraphs
#    (the occurrence inside the "数据来源" explanation paragraph is left as-is).
$p1 = $d.Paragraphs.Item(15)
$p1.Range.Find.Execute("localhost:8080/api/", $true, $false, $false, $false, $false, $true, 1, $false, "localhost:8090/api/", 2)

$p2 = $d.Paragraphs.Item(16)
$p2.Range.Find.Execute("localhost:8080/app/index.html", $true, $false, $false, $false, $false, $true, 1, $false, "localhost:8090/app/index.html", 2)

# 2. Move the "_GoBack" bookmark from the end of the document (after the
#    github link) to the blank paragraph right before "访问：" (the second
#    of the two empty "widowControl" paragraphs, immediately preceding it).
$target = $d.Paragraphs.Item(12).Range
$d.Bookmarks.Add("_GoBack", $target)
